$d = $word.ActiveDocument

# Manual line-break character used by Word for <w:br/> (Chr(11)).
$lb = [char]11

# Append a blank line (break-only run), matching the existing
# "break-only run" pattern already used earlier in the paragraph.
$r1 = $d.Content
$r1.Collapse(0)   # wdCollapseEnd
$r1.InsertAfter("" + $lb)
$r1.Collapse(0)

# Append a line break followed by "The app broke".
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertAfter("" + $lb + "The app broke")
$r2.Collapse(0)

# Append a line break followed by "we fixed the bug".
$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertAfter("" + $lb + "we fixed the bug")

$d.Save()
